$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SBASciFALL1819")
$ws.Activate()

# Update the test name / subtest values for the FALL sheet so they are
# distinguishable from the SPRING sheet's values.
$ws.Range("F19").Value = "SBASCI_FALL"
$ws.Range("F20").Value = "SCI_FALL"

# Remove the last data row (row 29, the "status" variable) entirely,
# shifting any rows below it upward (none exist here, so it just drops).
$ws.Rows.Item(29).Delete()

# Renumber the VarNum column (A) sequentially for the remaining rows.
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Update the active selection to match the new state of the sheet.
$ws.Range("A2:A28").Select()
